# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-14
$kValues = @{
    2  = 1
    3  = 3
    4  = 4
    5  = 5
    6  = 4
    7  = 4
    8  = 2
    9  = 7
    10 = 3
    11 = 4
    12 = 6
    13 = 5
    14 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
